# Updated cryptos list on Mon Nov 20 14:09:02 UTC 2023 with GitHub Actions
#
# The sheet stores Price (D) / Volume(1h) (E) figures as plain text cells
# (they include thousands-dot-separated numbers like "37.151.06" that are
# not valid numeric literals). When a value DOES look like a normal decimal
# number (e.g. "59.88"), Excel's usual text->value coercion would turn it
# into a Number cell on assignment, which would silently change the cell's
# type. To keep those cells as literal text (matching the source data) we
# prefix such values with a leading apostrophe, exactly like typing an
# apostrophe in the Excel UI to force text entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = [ordered]@{
    2  = @{ D = "37.151.06";  E = "  +1.71%  " }
    3  = @{ D = "2.018.02";   E = "  +3.38%  " }
    4  = @{ D = $null;        E = "  -0.08%  " }
    5  = @{ D = "246.19";     E = "  +1.31%  " }
    6  = @{ D = "0.625";      E = "  +0.10%  " }
    7  = @{ D = "59.88";      E = "  +0.14%  " }
    8  = @{ D = $null;        E = "  -0.03%  " }
    9  = @{ D = "0.391";      E = "  +3.60%  " }
    10 = @{ D = "0.0806";     E = "  +2.52%  " }
    11 = @{ D = $null;        E = "  +1.07%  " }
    12 = @{ D = "14.96";      E = "  +5.93%  " }
    13 = @{ D = "2.315.68";   E = "  +3.37%  " }
    14 = @{ D = "0.846";      E = "  +1.26%  " }
    15 = @{ D = "21.94";      E = "  +2.34%  " }
    16 = @{ D = "5.42";       E = "  +3.18%  " }
    17 = @{ D = "2.018.22";   E = "  +3.02%  " }
    18 = @{ D = "37.109.68";  E = "  +1.84%  " }
    19 = @{ D = "70.24";      E = "  +1.62%  " }
    20 = @{ D = "0.0₃0860"; E = "  +0.95%  " }
    21 = @{ D = "5.20";       E = "  +2.83%  " }
    22 = @{ D = "230.23";     E = "  +0.48%  " }
    23 = @{ D = "1.00";       E = "  +0.06%  " }
    24 = @{ D = "2.57";       E = "  +5.17%  " }
    25 = @{ D = "2.34";       E = "  -0.87%  " }
    26 = @{ D = "9.34";       E = "  +2.17%  " }
    27 = @{ D = "163.25";     E = "  +1.87%  " }
    28 = @{ D = "0.136";      E = "  -3.65%  " }
    29 = @{ D = "19.73";      E = "  +2.56%  " }
    30 = @{ D = "1.38";       E = "  +6.24%  " }
    31 = @{ D = $null;        E = "  +0.87%  " }
    32 = @{ D = "0.0670";     E = "  +9.78%  " }
    33 = @{ D = "4.75";       E = "  -0.02%  " }
    34 = @{ D = "2.51";       E = "  +11.55%  " }
    35 = @{ D = "4.44";       E = "  -0.09%  " }
    36 = @{ D = "3.58";       E = "  +5.63%  " }
    37 = @{ D = $null;        E = "  -0.18%  " }
    38 = @{ D = $null;        E = "  +1.57%  " }
    39 = @{ D = $null;        E = "  -1.90%  " }
    40 = @{ D = "3.00";       E = "  +3.10%  " }
    41 = @{ D = "0.0969";     E = "  +0.68%  " }
    42 = @{ D = $null;        E = "  +3.29%  " }
    43 = @{ D = $null;        E = "  +1.43%  " }
    44 = @{ D = "16.62";      E = "  +5.08%  " }
    45 = @{ D = "91.07";      E = "  +2.93%  " }
    46 = @{ D = "1.372.36";   E = "  +1.01%  " }
    47 = @{ D = $null;        E = "  +2.70%  " }
    48 = @{ D = $null;        E = "  +3.79%  " }
    49 = @{ D = "2.11";       E = "  +15.03%  " }
    50 = @{ D = $null;        E = "  +1.73%  " }
    51 = @{ D = "46.01";      E = "  -0.34%  " }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]

    if ($null -ne $vals.D) {
        $dText = $vals.D
        # If the raw text parses as a plain number, Excel would silently
        # convert the assigned string into a numeric cell. Force it to stay
        # text (quote-prefix) so the stored cell keeps its original text type
        # and exact formatting (e.g. "1.00" / "0.0670" keep trailing zeros).
        $looksNumeric = $dText -match "^[0-9]+(\.[0-9]+)?$"
        if ($looksNumeric) {
            $ws.Range("D$row").Value = "'" + $dText
        } else {
            $ws.Range("D$row").Value = $dText
        }
    }

    $ws.Range("E$row").Value = $vals.E
}
